$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this record block (rows 124-125),
# pushing the existing rows 124-145 down to 126-147.
$ws.Rows("124:125").Insert()

# New row 124: weekly price record (Especial) dated 2021-11-04 (serial 44504)
$ws.Cells.Item(124, 1).Value  = 3
$ws.Cells.Item(124, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(124, 3).Value  = "Coquimbo"
$ws.Cells.Item(124, 4).Value  = 44504
$ws.Cells.Item(124, 5).Value  = 5
$ws.Cells.Item(124, 6).Value  = "Fruta"
$ws.Cells.Item(124, 7).Value  = 100101
$ws.Cells.Item(124, 8).Value  = "Berries"
$ws.Cells.Item(124, 9).Value  = 100112025
$ws.Cells.Item(124, 10).Value = "Frutilla"
$ws.Cells.Item(124, 11).Value = "Sin especificar"
$ws.Cells.Item(124, 12).Value = "Especial"
$ws.Cells.Item(124, 13).Value = 125
$ws.Cells.Item(124, 14).Value = 6000
$ws.Cells.Item(124, 15).Value = 6000
$ws.Cells.Item(124, 16).Value = 6000
$ws.Cells.Item(124, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(124, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(124, 19).Value = 857
$ws.Cells.Item(124, 20).Value = 7

# New row 125: weekly price record (Segunda) dated 2021-11-04 (serial 44504)
$ws.Cells.Item(125, 1).Value  = 3
$ws.Cells.Item(125, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(125, 3).Value  = "Coquimbo"
$ws.Cells.Item(125, 4).Value  = 44504
$ws.Cells.Item(125, 5).Value  = 5
$ws.Cells.Item(125, 6).Value  = "Fruta"
$ws.Cells.Item(125, 7).Value  = 100101
$ws.Cells.Item(125, 8).Value  = "Berries"
$ws.Cells.Item(125, 9).Value  = 100112025
$ws.Cells.Item(125, 10).Value = "Frutilla"
$ws.Cells.Item(125, 11).Value = "Sin especificar"
$ws.Cells.Item(125, 12).Value = "Segunda"
$ws.Cells.Item(125, 13).Value = 75
$ws.Cells.Item(125, 14).Value = 4000
$ws.Cells.Item(125, 15).Value = 4000
$ws.Cells.Item(125, 16).Value = 4000
$ws.Cells.Item(125, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(125, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(125, 19).Value = 571
$ws.Cells.Item(125, 20).Value = 7
